$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Update Date value (row 8)
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row after "Contact" (row 10) for the new "Jurisdiction" property
$ws.Rows.Item(11).Insert()

# Populate the new row
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Copy formatting from the row above so the new row matches the existing style
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false
